$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

$rng = $ws.Range("A3:A12")
$rng.NumberFormat = "@"
$rng.Value = "16"
$rng.Style = "Normal"
